# Update ACTUAL_ENERGY (column B) and recompute WESM_EXPOSURE (column D)
# for hours 20-23 (rows 21-24), as part of "Added DAP for today".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ACTUAL_ENERGY values
$ws.Range("B21").Value = 24833.09097147339
$ws.Range("B22").Value = 26245.15567101265
$ws.Range("B23").Value = 29155.857
$ws.Range("B24").Value = 31722.2145

# Recomputed WESM_EXPOSURE = ACTUAL_ENERGY - TOTAL_BCQ_NOMINATION
$ws.Range("D21").Value = -50166.90902852661
$ws.Range("D22").Value = -48754.84432898735
$ws.Range("D23").Value = -45844.143
$ws.Range("D24").Value = -40277.7855
